$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data occupies rows 2-8 (after the header row 1), columns A:D.
# The edit re-sorts the data rows in ascending order by column A (time).
$dataRange = $ws.Range("A2:D8")
$sortField = $ws.Range("A2:A8")

$dataRange.Sort(
    $sortField,            # Key1
    1,                     # Order1 = xlAscending
    [System.Type]::Missing,# Key2
    [System.Type]::Missing,# Type
    [System.Type]::Missing,# Order2
    [System.Type]::Missing,# Key3
    [System.Type]::Missing,# Order3
    2                       # Header = xlNo
)
